$d = $word.ActiveDocument

# 1. Title change: "Technology of Biopolymers" -> "Biopolymers Technology"
$d.Content.Find.Execute("Technology of Biopolymers", $true, $false, $false, $false, $false, $true, 1, $false, "Biopolymers Technology", 2) | Out-Null

# 2. Activation date change
$d.Content.Find.Execute("Ativação: 01/01/2017", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2025", 2) | Out-Null

# Helper function to insert a new italic paragraph right after a given paragraph,
# matching the target OOXML shape: <w:p><w:r><w:rPr><w:i/></w:rPr><w:t>...</w:t></w:r></w:p>
function Insert-ItalicParagraphAfter($para, $text) {
    $para.Range.InsertParagraphAfter()
    $newPara = $para.Next()
    $newRange = $newPara.Range
    $endChar = $newRange.End
    $textRange = $d.Range($newRange.Start, $endChar - 1)
    $textRange.Text = $text
    $textRange.Font.Italic = $true
}

# 3. Objetivos paragraph: replace text, then add new italic English paragraph after it
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Apresentar os conceitos básicos da ciência dos polímeros, incluindo*") {
        $p.Range.Text = "Apresentar os conceitos básicos da ciência dos polímeros e os principais problemas diretamente relacionados ao seu uso e descarte indiscriminados. Desenvolver o pensamento crítico e apresentar ferramentas alternativas para a produção dos polímeros, bem como para minimizar seus impactos ambientais."
        Insert-ItalicParagraphAfter $p "To present the basic concepts of polymer science and the main problems directly related to their indiscriminate use and disposal. To develop critical thinking and present alternative tools for polymer production, as well as to minimize their environmental impacts."
        break
    }
}

# 4. Programa resumido paragraph: replace text, then add new italic English paragraph after it
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Fundamentos sobre a química dos polímeros; Mecanismos de polimerização*") {
        $p.Range.Text = "Reações de polimerização; Propriedades gerais dos polímeros; Monômeros e polímeros derivados de fontes renováveis. Conceitos básicos de circularidade dos materiais poliméricos."
        Insert-ItalicParagraphAfter $p "Polymerization reactions; General properties of polymers; Monomers and polymers derived from renewable sources. Basic concepts of circularity of polymeric materials."
        break
    }
}

# 5. Programa (full) paragraph: replace text, then add new italic English paragraph after it
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Fundamentos sobre a química dos polímeros: composição e estrutura, nomenclatura*") {
        $p.Range.Text = "Fundamentos sobre a química dos polímeros: composição e estrutura, massa molecular média, propriedades físicas (comportamentos cristalino e amorfo, propriedades mecânicas e térmicas). Reações de poliadição e policondensação. Introdução aos materiais derivados de fontes renováveis. Monômeros de fonte renovável (etileno, ácidos carboxílicos, aminas, álcoois, óleos vegetais, CO2, entre outros). Polímeros de fonte renovável (celulose, amido, quitina e quitosana, exopolissacarídeos, polihidroxialcanoatos). Estratégias para fim de vida: conceitos básicos de biodegradação e economia circular."
        Insert-ItalicParagraphAfter $p "Fundamentals of polymer chemistry: composition and structure, average molecular weight, physical properties (crystalline and amorphous behaviours, mechanical and thermal properties). Polyaddition and polycondensation reactions. Introduction to materials derived from renewable sources. Renewable source monomers (ethylene, carboxylic acids, amines, alcohols, vegetable oils, CO2, among others). Polymers directly extracted from renewable sources (cellulose, starch, chitin and chitosan, exopolysaccharides, polyhydroxyalkanoates). End-of-life strategies: basic concepts of biodegradation and circular economy."
        break
    }
}

# 6. Avaliação - Método text
$d.Content.Find.Execute("Duas provas escritas envolvendo o conteúdo teórico ministrado em sala de aula.", $true, $false, $false, $false, $false, $true, 1, $false, "Uma avaliação escrita e um estudo de caso.", 2) | Out-Null

# 7. Avaliação - Critério text
$d.Content.Find.Execute("A nota final corresponderá à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados.", $true, $false, $false, $false, $false, $true, 1, $false, "A nota final corresponderá à média aritmética da nota da prova escrita e da nota do estudo de caso. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto aqueles que tiverem média inferior a 3 estarão reprovados.", 2) | Out-Null

# 8. Bibliografia paragraph replacement
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Eloisa B. Mano*") {
        $p.Range.Text = "1-Eloisa B. Mano; Introdução a Polímeros, Editora Edgard BlücherLtda, São Paulo, 1999. 2-2- Sebastião V. Canevarol; Ciência dos Polímeros. Um Texto Básico Para Tecnólogos e Engenheiros. Artliber; 3ª edição. 3-3- J. P. Greene; Sustainable plastics: environmental assessments of biobased, biodegradable, and reclycled plastics. John Wiley & Sons, New Jersey, United States, 2014."
        break
    }
}
